$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the "K" column quote-prefixed "=" labels to the existing
#     AND (rows 4-6) and OR (rows 10-13) truth tables, mirroring the
#     style already used on row 3 / row 10's K column. ---
$ws.Range("K4").Value = "'="
$ws.Range("K5").Value = "'="
$ws.Range("K6").Value = "'="

$ws.Range("K10").Value = "'="
$ws.Range("K11").Value = "'="
$ws.Range("K12").Value = "'="
$ws.Range("K13").Value = "'="

# --- New NOT perceptron block (rows 15-19), mirroring the layout of
#     the AND block above it. ---

# Row 15: truth-table row + weight header values
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "NOT"
$ws.Range("E15").Value = "'="
$ws.Range("F15").Value = 0
$ws.Range("H15").Value = -30
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 15

# Row 16: second truth-table row + first computation row
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = "NOT"
$ws.Range("E16").Value = "'="
$ws.Range("F16").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("K16").Value = "'="
$ws.Range("L16").Value = -15
$ws.Range("N16").Value = 0

# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 1
$ws.Range("K17").Value = "'="
$ws.Range("L17").Value = 15
$ws.Range("N17").Value = 1

# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = "'="
$ws.Range("L18").Value = 15
$ws.Range("N18").Value = 1

# Row 19
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = "'="
$ws.Range("L19").Value = -15
$ws.Range("N19").Value = 0

# --- Cosmetic follow-ups that mirror what Excel itself would do after
#     this edit: widen column C (bestFit) so the new "NOT" label isn't
#     clipped, and leave the selection where the user's last edit was. ---
$ws.Columns("C").ColumnWidth = 4

$null = $ws.Range("N20").Select()
